$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 291, shifting rows 291:321 down to 292:322
$ws.Rows("291:291").Insert()

# Populate the newly inserted row 291 with the new record
$ws.Range("A291").Value = 11
$ws.Range("B291").Value = 'Vega Monumental Concepción'
$ws.Range("C291").Value = 'Bíobío'
$ws.Range("D291").Value = Get-Date -Year 2022 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("E291").Value = 8
$ws.Range("F291").Value = 'Fruta'
$ws.Range("G291").Value = 100102
$ws.Range("H291").Value = 'Cítricos'
$ws.Range("I291").Value = 100102005
$ws.Range("J291").Value = 'Naranja'
$ws.Range("K291").Value = 'Lane Late'
$ws.Range("L291").Value = 'Primera'
$ws.Range("M291").Value = 300
$ws.Range("N291").Value = 5500
$ws.Range("O291").Value = 6000
$ws.Range("P291").Value = 5750
$ws.Range("Q291").Value = '$/bandeja 15 kilos granel'
$ws.Range("R291").Value = "Región de O'Higgins"
$ws.Range("S291").Value = 383
$ws.Range("T291").Value = 15
